# Auto-generated Excel COM-interop edit script
# Updates crypto price/volume(1h) figures and fixes the WEMIXTOKEN/FraxShare row order
# (GitHub Actions scheduled refresh of cryptos.xlsx)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
if ($ws -eq $null) { $ws = $wb.Worksheets.Item(1) }

# Helper: write a literal TEXT value into a cell without Excel re-typing it
# as a number (keeps "24.640.59"-style and "  +0.28%  "-style strings intact,
# and leaves the cells style/number-format completely untouched).
function Set-TextValue {
    param($addr, $text)
    $escaped = $text -replace '"', '""'
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4163) | Out-Null
}

$excel.CutCopyMode = $false

Set-TextValue "D2" "24.640.59"
Set-TextValue "E2" "  +0.28%  "
Set-TextValue "D3" "1.697.00"
Set-TextValue "E3" "  +0.22%  "
Set-TextValue "E4" "  +0.12%  "
Set-TextValue "D5" "315.42"
Set-TextValue "E5" "  -0.38%  "
Set-TextValue "D6" "1.001"
Set-TextValue "E6" "  +0.01%  "
Set-TextValue "D7" "0.3924"
Set-TextValue "E7" "  -0.15%  "
Set-TextValue "D8" "0.4033"
Set-TextValue "E8" "  +0.47%  "
Set-TextValue "E9" "  -0.78%  "
Set-TextValue "E10" "  -0.08%  "
Set-TextValue "D11" "53.08"
Set-TextValue "E11" "  -1.73%  "
Set-TextValue "D12" "0.08833"
Set-TextValue "E12" "  +0.91%  "
Set-TextValue "D13" "7.448"
Set-TextValue "E13" "  +3.31%  "
Set-TextValue "D14" "23.60"
Set-TextValue "E14" "  +1.66%  "
Set-TextValue "D15" "8.224"
Set-TextValue "E15" "  +8.08%  "
Set-TextValue "D16" "0.00001318"
Set-TextValue "E16" "  -0.16%  "
Set-TextValue "D17" "1.705.53"
Set-TextValue "E17" "  +0.43%  "
Set-TextValue "D18" "99.55"
Set-TextValue "E18" "  -1.18%  "
Set-TextValue "D20" "19.64"
Set-TextValue "E20" "  +0.18%  "
Set-TextValue "D21" "7.081"
Set-TextValue "E21" "  +3.33%  "
Set-TextValue "E22" "  +0.39%  "
Set-TextValue "D23" "14.62"
Set-TextValue "E23" "  +4.18%  "
Set-TextValue "D24" "24.663.03"
Set-TextValue "E24" "  +0.37%  "
Set-TextValue "D25" "3.113"
Set-TextValue "E25" "  +3.17%  "
Set-TextValue "D26" "2.357"
Set-TextValue "E26" "  +1.80%  "
Set-TextValue "D27" "22.59"
Set-TextValue "E27" "  +1.04%  "
Set-TextValue "D28" "162.76"
Set-TextValue "E28" "  +2.27%  "
Set-TextValue "D29" "8.641"
Set-TextValue "E29" "  +15.09%  "
Set-TextValue "D30" "135.56"
Set-TextValue "E30" "  +1.36%  "
Set-TextValue "D31" "5.145"
Set-TextValue "E31" "  -1.11%  "
Set-TextValue "D32" "0.08952"
Set-TextValue "E32" "  +5.07%  "
Set-TextValue "D33" "7.571"
Set-TextValue "E33" "  +3.56%  "
Set-TextValue "D34" "1.065"
Set-TextValue "E34" "  -2.89%  "
Set-TextValue "B35" "FraxShare"
Set-TextValue "C35" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D35" "11.07"
Set-TextValue "E35" "  -2.31%  "
Set-TextValue "B36" "WEMIXTOKEN"
Set-TextValue "C36" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D36" "1.963"
Set-TextValue "E36" "  -0.60%  "
Set-TextValue "D37" "0.2753"
Set-TextValue "E37" "  +1.20%  "
Set-TextValue "D38" "14.42"
Set-TextValue "E38" "  -0.79%  "
Set-TextValue "D39" "0.02819"
Set-TextValue "E39" "  +2.77%  "
Set-TextValue "D40" "0.09113"
Set-TextValue "E40" "  +0.96%  "
Set-TextValue "D41" "1.459"
Set-TextValue "E41" "  -0.57%  "
Set-TextValue "D42" "0.7647"
Set-TextValue "E42" "  -0.33%  "
Set-TextValue "D43" "15.84"
Set-TextValue "E43" "  +3.21%  "
Set-TextValue "D44" "0.7152"
Set-TextValue "E44" "  -0.41%  "
Set-TextValue "D45" "2.548"
Set-TextValue "E45" "  +1.77%  "
Set-TextValue "D46" "4.212"
Set-TextValue "E46" "  -0.08%  "
Set-TextValue "D47" "1.002"
Set-TextValue "E47" "  +0.06%  "
Set-TextValue "D48" "1.345"
Set-TextValue "E48" "  -0.13%  "
Set-TextValue "D49" "139.72"
Set-TextValue "E49" "  -1.09%  "
Set-TextValue "E50" "  -0.73%  "
Set-TextValue "D51" "90.20"
Set-TextValue "E51" "  +2.16%  "

$excel.CutCopyMode = $false
